# "Generate Report for Archive"
# The handoff status moved on: cells that used to read "Ready for handoff"
# now report "In Translation", and the Status column(s) are re-sized to fit
# the new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- Update status text on every sheet (Overview, zh-cn, de-de) ---------
foreach ($ws in $wb.Worksheets) {
    $found = $ws.Cells.Find("Ready for handoff")
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        do {
            $found.Value = "In Translation"
            $found = $ws.Cells.FindNext($found)
        } while ($found -ne $null -and $found.Address() -ne $firstAddress)
    }
}

# --- Re-fit the Status column(s) now that the text is shorter -----------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
